# Lab5 slide "Classes" (slide 4): fill in the previously-empty content
# placeholder with the "Three Classes" / "Information Hidden" bullet lists.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(1)     # "Content Placeholder 1"
$tr = $shape.TextFrame.TextRange

$ndash = [char]0x2013          # "–" used in "Activator – ..." / "Driver – ..."

# Build the text paragraph by paragraph. Starting the range's Text and then
# repeatedly calling InsertAfter appends a new run each time; a leading "`r"
# starts a brand-new paragraph, while text with no leading "`r" is appended
# as a second run within the current paragraph (used for the two paragraphs
# that mix a "clean" run with an "err"-flagged / plain run in the source deck).
$tr.Text = "Three Classes"
[void]$tr.InsertAfter("`rActivator $ndash Main program, deals with connection and timer")
[void]$tr.InsertAfter("`rDriver $ndash controls Robot movement and access to hardware")
[void]$tr.InsertAfter("`rMessageHandler")
[void]$tr.InsertAfter(" $ndash deals with decoding and creating messages")
[void]$tr.InsertAfter("`rInformation Hidden")
[void]$tr.InsertAfter("`rHow the robot is being controlled")
[void]$tr.InsertAfter("`rHow messages are created and decoded, and what they mean in terms of ")
[void]$tr.InsertAfter("robot movement")

# Second-level (sub-bullet) paragraphs.
$tr.Paragraphs(2,1).IndentLevel = 2   # Activator ...
$tr.Paragraphs(3,1).IndentLevel = 2   # Driver ...
$tr.Paragraphs(4,1).IndentLevel = 2   # MessageHandler ...
$tr.Paragraphs(6,1).IndentLevel = 2   # How the robot is being controlled
$tr.Paragraphs(7,1).IndentLevel = 2   # How messages are created and decoded ...
